# Moved DocumentDB and Storage to the HA/Common resource group deployment.
# On the architecture slide, the "Resource Group" textbox (id 114) is pulled
# into the "Group 12" group (id 13) -- right after the "Rectangle 14" shape
# -- and the group is shifted/resized to the new tight bounding box. The
# four elbow connectors (163, 166, 174, 177) that used to sit beside the
# group move to become trailing siblings right after the (re-created) group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# 1) Move the standalone "Resource Group" textbox (id 114) to its new
#    slide-absolute position before it gets folded into the group -- this
#    makes the group's auto-computed bounding box land exactly on the
#    target off/ext/chOff/chExt once everything is (re)grouped below.
$tb = $s.Shapes.Item("TextBox 113")
$tb.Left = 3366841 / 12700.0
$tb.Top = 5670658 / 12700.0

# 2) Ungroup "Group 12" (id 13) so its members become loose top-level
#    shapes again, in their original relative (document) order.
$grp = $s.Shapes.Item("Group 12")
$null = $grp.Ungroup()

# 3) Reposition the textbox in z-order so it sits immediately after
#    "Rectangle 14" -- matching the target child order inside the group.
$tb.ZOrder(0)
for ($i = 1; $i -le 8; $i++) {
    $tb.ZOrder(3)
}

# 4) Re-group "Rectangle 14" + the "Resource Group" textbox + the other
#    8 former group members. PowerPoint hands the freshly-made group the
#    smallest free shape id (2) and default name ("Group 1"), exactly
#    matching the target.
$names = @(
    "Rectangle 14",
    "TextBox 113",
    "Picture 98",
    "TextBox 101",
    "Picture 130",
    "TextBox 131",
    "Picture 107",
    "TextBox 111",
    "Picture 9",
    "TextBox 112"
)
$range = $s.Shapes.Range($names)
$newGroup = $range.Group()

# 5) Move the four elbow connectors so they become trailing top-level
#    siblings right after the regrouped shape, in the target order.
$s.Shapes.Item("Connector: Elbow 165").ZOrder(0)
$s.Shapes.Item("Connector: Elbow 162").ZOrder(0)
$s.Shapes.Item("Connector: Elbow 173").ZOrder(0)
$s.Shapes.Item("Connector: Elbow 176").ZOrder(0)
